$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 64. This shifts the existing
# rows 64-69 down to 65-70, preserving all of their data and formatting
# (matches the diff: old row 64 data now appears at row 65, ..., old row
# 69 data now appears at row 70).
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Range("A64").Value = 9
$ws.Range("B64").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C64").Value = "Metropolitana"
$ws.Range("D64").Value = 44714
$ws.Range("E64").Value = 13
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100102
$ws.Range("H64").Value = "Cítricos"
$ws.Range("I64").Value = 100102006
$ws.Range("J64").Value = "Pomelo"
$ws.Range("K64").Value = "Start Ruby"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 470
$ws.Range("N64").Value = 7500
$ws.Range("O64").Value = 8000
$ws.Range("P64").Value = 7766
$ws.Range("Q64").Value = "$/caja 14 kilos"
$ws.Range("R64").Value = "Región Metropolitana"
$ws.Range("S64").Value = 555
$ws.Range("T64").Value = 14
